$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 6 (old rows 6-9 shift down to 8-11).
$ws.Rows("6:7").Insert()

# New row 6: Primera quality entry for the 2022-01-13 (serial 44574) price date.
$row6 = @(11, "Vega Monumental Concepción", "Bíobío", 44574, 8, "Fruta", 100101, "Berries", 100101004, "Frambuesa", "Sin especificar", "Primera", 200, 7000, 8000, 7500, "`$/bandeja 2 kilos", "Región de Ñuble", 3750, 2)
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(6, $col).Value = $row6[$col - 1]
}

# New row 7: Segunda quality entry for the same price date.
$row7 = @(11, "Vega Monumental Concepción", "Bíobío", 44574, 8, "Fruta", 100101, "Berries", 100101004, "Frambuesa", "Sin especificar", "Segunda", 100, 6000, 6000, 6000, "`$/bandeja 2 kilos", "Región de Ñuble", 3000, 2)
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(7, $col).Value = $row7[$col - 1]
}
